$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Team Month")

# Update Repayment (column D) values for rows 2-8
$ws.Range("D2").Value = 760238652
$ws.Range("D3").Value = 307674510
$ws.Range("D4").Value = 200375241
$ws.Range("D5").Value = 187022090
$ws.Range("D6").Value = 362590491
$ws.Range("D7").Value = 656856222
$ws.Range("D8").Value = 637943333

# Update the active cell selection to D2
$ws.Range("D2").Select()
